$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update weekly HIGH/LOW/CLOSE/LTP/VOL/9:25 CLOSE data (rows 2-30)
$ws.Cells.Item(2, 2).Value = 2232.55
$ws.Cells.Item(2, 3).Value = 2185.05
$ws.Cells.Item(2, 4).Value = 2197.6
$ws.Cells.Item(2, 5).Value = 2196.65
$ws.Cells.Item(2, 6).Value = 18
$ws.Cells.Item(2, 7).Value = 2188.1
$ws.Cells.Item(3, 2).Value = 428.6
$ws.Cells.Item(3, 3).Value = 423.75
$ws.Cells.Item(3, 4).Value = 424
$ws.Cells.Item(3, 5).Value = 424.6
$ws.Cells.Item(3, 6).Value = 8
$ws.Cells.Item(3, 7).Value = 426.3
$ws.Cells.Item(4, 2).Value = 1616.45
$ws.Cells.Item(4, 3).Value = 1594.5
$ws.Cells.Item(4, 4).Value = 1607
$ws.Cells.Item(4, 5).Value = 1610.4
$ws.Cells.Item(4, 6).Value = 17
$ws.Cells.Item(4, 7).Value = 1604.25
$ws.Cells.Item(5, 2).Value = 7128.85
$ws.Cells.Item(5, 3).Value = 7049
$ws.Cells.Item(5, 4).Value = 7100
$ws.Cells.Item(5, 5).Value = 7103.15
$ws.Cells.Item(5, 6).Value = 16
$ws.Cells.Item(5, 7).Value = 7057.5
$ws.Cells.Item(6, 2).Value = 215.2
$ws.Cells.Item(6, 3).Value = 213.25
$ws.Cells.Item(6, 4).Value = 213.8
$ws.Cells.Item(6, 5).Value = 213.9
$ws.Cells.Item(6, 6).Value = 56
$ws.Cells.Item(6, 7).Value = 215.1
$ws.Cells.Item(7, 2).Value = 198.45
$ws.Cells.Item(7, 3).Value = 194.65
$ws.Cells.Item(7, 4).Value = 195.5
$ws.Cells.Item(7, 5).Value = 195.75
$ws.Cells.Item(7, 6).Value = 89
$ws.Cells.Item(7, 7).Value = 198.05
$ws.Cells.Item(8, 2).Value = 338.5
$ws.Cells.Item(8, 3).Value = 328.85
$ws.Cells.Item(8, 4).Value = 334.4
$ws.Cells.Item(8, 5).Value = 334.2
$ws.Cells.Item(8, 6).Value = 151
$ws.Cells.Item(8, 7).Value = 335.4
$ws.Cells.Item(9, 2).Value = 641.2
$ws.Cells.Item(9, 3).Value = 627.15
$ws.Cells.Item(9, 4).Value = 633.35
$ws.Cells.Item(9, 5).Value = 634.45
$ws.Cells.Item(9, 6).Value = 60
$ws.Cells.Item(9, 7).Value = 631.65
$ws.Cells.Item(10, 2).Value = 3863.8
$ws.Cells.Item(10, 3).Value = 3826.05
$ws.Cells.Item(10, 4).Value = 3835
$ws.Cells.Item(10, 5).Value = 3844.2
$ws.Cells.Item(10, 6).Value = 2
$ws.Cells.Item(10, 7).Value = 3855
$ws.Cells.Item(11, 2).Value = 148.2
$ws.Cells.Item(11, 3).Value = 146.5
$ws.Cells.Item(11, 4).Value = 146.55
$ws.Cells.Item(11, 5).Value = 146.7
$ws.Cells.Item(11, 6).Value = 77
$ws.Cells.Item(11, 7).Value = 147.55
$ws.Cells.Item(12, 2).Value = 1334.4
$ws.Cells.Item(12, 3).Value = 1321.8
$ws.Cells.Item(12, 4).Value = 1326.05
$ws.Cells.Item(12, 5).Value = 1326.6
$ws.Cells.Item(12, 6).Value = 15
$ws.Cells.Item(12, 7).Value = 1331.5
$ws.Cells.Item(13, 2).Value = 1522
$ws.Cells.Item(13, 3).Value = 1513.4
$ws.Cells.Item(13, 4).Value = 1518.25
$ws.Cells.Item(13, 5).Value = 1517.95
$ws.Cells.Item(13, 6).Value = 124
$ws.Cells.Item(13, 7).Value = 1516.75
$ws.Cells.Item(14, 2).Value = 510.95
$ws.Cells.Item(14, 3).Value = 505.45
$ws.Cells.Item(14, 4).Value = 506.55
$ws.Cells.Item(14, 5).Value = 506.85
$ws.Cells.Item(14, 6).Value = 71
$ws.Cells.Item(14, 7).Value = 505.95
$ws.Cells.Item(15, 2).Value = 928.05
$ws.Cells.Item(15, 3).Value = 922
$ws.Cells.Item(15, 4).Value = 927.45
$ws.Cells.Item(15, 5).Value = 926.1
$ws.Cells.Item(15, 6).Value = 89
$ws.Cells.Item(15, 7).Value = 925.05
$ws.Cells.Item(16, 2).Value = 1507
$ws.Cells.Item(16, 3).Value = 1494.1
$ws.Cells.Item(16, 4).Value = 1502
$ws.Cells.Item(16, 5).Value = 1502.05
$ws.Cells.Item(16, 6).Value = 16
$ws.Cells.Item(16, 7).Value = 1506.15
$ws.Cells.Item(17, 2).Value = 1448.65
$ws.Cells.Item(17, 3).Value = 1436.75
$ws.Cells.Item(17, 4).Value = 1439.4
$ws.Cells.Item(17, 5).Value = 1439
$ws.Cells.Item(17, 6).Value = 37
$ws.Cells.Item(17, 7).Value = 1446.25
$ws.Cells.Item(18, 2).Value = 660.5
$ws.Cells.Item(18, 3).Value = 652.7
$ws.Cells.Item(18, 4).Value = 654.35
$ws.Cells.Item(18, 5).Value = 654.65
$ws.Cells.Item(18, 6).Value = 12
$ws.Cells.Item(18, 7).Value = 656
$ws.Cells.Item(19, 2).Value = 465.8
$ws.Cells.Item(19, 3).Value = 460.5
$ws.Cells.Item(19, 4).Value = 461.8
$ws.Cells.Item(19, 5).Value = 462.05
$ws.Cells.Item(19, 6).Value = 8
$ws.Cells.Item(19, 7).Value = 464.8
$ws.Cells.Item(20, 2).Value = 1571.4
$ws.Cells.Item(20, 3).Value = 1551.2
$ws.Cells.Item(20, 4).Value = 1561.25
$ws.Cells.Item(20, 5).Value = 1559.75
$ws.Cells.Item(20, 6).Value = 22
$ws.Cells.Item(20, 7).Value = 1556.4
$ws.Cells.Item(21, 2).Value = 276.5
$ws.Cells.Item(21, 3).Value = 272.4
$ws.Cells.Item(21, 4).Value = 275.15
$ws.Cells.Item(21, 5).Value = 275.15
$ws.Cells.Item(21, 6).Value = 20
$ws.Cells.Item(21, 7).Value = 275.4
$ws.Cells.Item(22, 2).Value = 251.65
$ws.Cells.Item(22, 3).Value = 247.3
$ws.Cells.Item(22, 4).Value = 249.75
$ws.Cells.Item(22, 5).Value = 250
$ws.Cells.Item(22, 6).Value = 68
$ws.Cells.Item(22, 7).Value = 250.55
$ws.Cells.Item(23, 2).Value = 2388
$ws.Cells.Item(23, 3).Value = 2364.1
$ws.Cells.Item(23, 4).Value = 2377
$ws.Cells.Item(23, 5).Value = 2378.9
$ws.Cells.Item(23, 6).Value = 41
$ws.Cells.Item(23, 7).Value = 2365.55
$ws.Cells.Item(24, 2).Value = 566.2
$ws.Cells.Item(24, 3).Value = 561
$ws.Cells.Item(24, 4).Value = 561.4
$ws.Cells.Item(24, 5).Value = 561.5
$ws.Cells.Item(24, 6).Value = 142
$ws.Cells.Item(24, 7).Value = 565.75
$ws.Cells.Item(25, 2).Value = 676.7
$ws.Cells.Item(25, 3).Value = 660.5
$ws.Cells.Item(25, 4).Value = 663.1
$ws.Cells.Item(25, 5).Value = 662.55
$ws.Cells.Item(25, 6).Value = 4
$ws.Cells.Item(25, 7).Value = 675
$ws.Cells.Item(26, 2).Value = 971
$ws.Cells.Item(26, 3).Value = 960.15
$ws.Cells.Item(26, 4).Value = 963.3
$ws.Cells.Item(26, 5).Value = 961.3
$ws.Cells.Item(26, 6).Value = 3
$ws.Cells.Item(26, 7).Value = 962
$ws.Cells.Item(27, 2).Value = 683.3
$ws.Cells.Item(27, 3).Value = 675.25
$ws.Cells.Item(27, 4).Value = 681
$ws.Cells.Item(27, 5).Value = 681.7
$ws.Cells.Item(27, 6).Value = 72
$ws.Cells.Item(27, 7).Value = 677.4
$ws.Cells.Item(28, 2).Value = 263.85
$ws.Cells.Item(28, 3).Value = 260.75
$ws.Cells.Item(28, 4).Value = 262.6
$ws.Cells.Item(28, 5).Value = 262.65
$ws.Cells.Item(28, 6).Value = 69
$ws.Cells.Item(28, 7).Value = 263.55
$ws.Cells.Item(29, 2).Value = 126.4
$ws.Cells.Item(29, 3).Value = 125.55
$ws.Cells.Item(29, 4).Value = 126.2
$ws.Cells.Item(29, 5).Value = 126.25
$ws.Cells.Item(29, 6).Value = 296
$ws.Cells.Item(29, 7).Value = 125.75
$ws.Cells.Item(30, 2).Value = 8724.35
$ws.Cells.Item(30, 3).Value = 8665
$ws.Cells.Item(30, 4).Value = 8704
$ws.Cells.Item(30, 5).Value = 8708.35
$ws.Cells.Item(30, 7).Value = 8717.950000000001

# Update the active cell selection on Sheet1
$ws.Range("I7").Select() | Out-Null
